# Auto-generated cell value updates for cryptos.xlsx snapshot refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole data block to Text format first so numeric-looking
# strings (prices like "1.012", percentages, tiny decimals) are written
# back as literal text instead of being auto-coerced into numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.688.05"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.016.01"
$ws.Range("E3").Value = "  -4.72%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "331.77"
$ws.Range("E5").Value = "  -4.36%  "
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "0.4999"
$ws.Range("E7").Value = "  -3.66%  "
$ws.Range("D8").Value = "0.4204"
$ws.Range("E8").Value = "  -5.45%  "
$ws.Range("D9").Value = "54.33"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "0.08910"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").Value = "1.126"
$ws.Range("E11").Value = "  -4.69%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.311.50"
$ws.Range("E12").Value = "  +8.37%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "23.18"
$ws.Range("E13").Value = "  -7.84%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "8.178"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "6.532"
$ws.Range("E15").Value = "  -5.32%  "
$ws.Range("D16").Value = "97.23"
$ws.Range("E16").Value = "  -5.39%  "
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "0.00001113"
$ws.Range("E18").Value = "  -4.12%  "
$ws.Range("D19").Value = "0.06619"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "19.57"
$ws.Range("E20").Value = "  -9.00%  "
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "6.027"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").Value = "29.698.25"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "11.93"
$ws.Range("E24").Value = "  -5.92%  "
$ws.Range("D25").Value = "2.308"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "2.233.79"
$ws.Range("E26").Value = "  -6.25%  "
$ws.Range("D27").Value = "158.32"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "20.75"
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("D29").Value = "6.549"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "2.344"
$ws.Range("E30").Value = "  -7.64%  "
$ws.Range("D31").Value = "127.78"
$ws.Range("E31").Value = "  -4.81%  "
$ws.Range("D32").Value = "1.062"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("D33").Value = "0.09946"
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("D34").Value = "1.560"
$ws.Range("E34").Value = "  -11.70%  "
$ws.Range("D35").Value = "3.860"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").Value = "5.813"
$ws.Range("E36").Value = "  -6.93%  "
$ws.Range("D37").Value = "9.715"
$ws.Range("E37").Value = "  -9.47%  "
$ws.Range("D38").Value = "0.02466"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").Value = "0.06404"
$ws.Range("E39").Value = "  -6.74%  "
$ws.Range("D40").Value = "1.292"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").Value = "11.89"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("D42").Value = "0.6552"
$ws.Range("E42").Value = "  -7.28%  "
$ws.Range("D43").Value = "0.2084"
$ws.Range("E43").Value = "  -7.13%  "
$ws.Range("D44").Value = "1.012"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "0.6362"
$ws.Range("E45").Value = "  -7.06%  "
$ws.Range("D46").Value = "2.239"
$ws.Range("E46").Value = "  -5.10%  "
$ws.Range("D47").Value = "13.38"
$ws.Range("E47").Value = "  -7.82%  "
$ws.Range("D48").Value = "1.276"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "3.570"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07047"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.00000000327"
$ws.Range("E51").Value = "  -8.47%  "

# Restore the default (un-styled) cell style now that the text values are
# committed, so we do not leave a stray numeric format applied to cells.
$dataRange.Style = "Normal"

